$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1949.9166
$ws.Range("J17").Value = 1949.9166
$ws.Range("L17").Value = 5849.7498
$ws.Range("N17").Value = -6185.7498
$ws.Range("H64").Value = 7191.3
$ws.Range("I64").Value = 5534.8
$ws.Range("K64").Value = 5534.8
$ws.Range("M64").Value = -5286.8
$ws.Range("H67").Value = 7191.3
$ws.Range("I67").Value = 5534.8
$ws.Range("K67").Value = 5534.8
$ws.Range("M67").Value = -4676.8
$ws.Range("H69").Value = 25394.055
$ws.Range("H72").Value = 25394.055
$ws.Range("H86").Value = 4390439.5
$ws.Range("J86").Value = 7523225
$ws.Range("L86").Value = 7523225
$ws.Range("N86").Value = -7525471
$ws.Range("H89").Value = 4390439.5
$ws.Range("J89").Value = 7523225
$ws.Range("L89").Value = 37616125
$ws.Range("N89").Value = -37627357
$ws.Range("H92").Value = 4808525
$ws.Range("I92").Value = 638.05
$ws.Range("J92").Value = 20834816
$ws.Range("K92").Value = 638.05
$ws.Range("L92").Value = 20834816
$ws.Range("M92").Value = 609.95
$ws.Range("N92").Value = -20837312
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
$ws.Range("H115").Value = 1297.8572
$ws.Range("I115").Value = 1297.8572
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3893.5716
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -2326.5716
$ws.Range("N115").ClearContents()
$ws.Range("H135").Value = 7052.4443
$ws.Range("I135").Value = 2985.0667
$ws.Range("K135").Value = 26865.6003
$ws.Range("M135").Value = -24330.6003
$ws.Range("H137").Value = 4244.759
$ws.Range("I137").Value = 7977.5
$ws.Range("K137").Value = 23932.5
$ws.Range("M137").Value = -21382.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9285.405000000001
$ws.Range("I32").Value = 9285.405000000001
$ws.Range("K32").Value = 9285.405000000001
$ws.Range("M32").Value = -8998.405000000001
$ws.Range("H61").Value = 5867.3335
$ws.Range("I61").Value = 4344
$ws.Range("K61").Value = 4344
$ws.Range("M61").Value = -4132
$ws.Range("H63").Value = 6097
$ws.Range("I63").Value = 9170
$ws.Range("J63").Value = 1487.5
$ws.Range("K63").Value = 9170
$ws.Range("L63").Value = 1487.5
$ws.Range("M63").Value = -8484
$ws.Range("N63").Value = -2859.5
$ws.Range("H66").Value = 6097
$ws.Range("I66").Value = 9170
$ws.Range("J66").Value = 1487.5
$ws.Range("K66").Value = 45850
$ws.Range("L66").Value = 7437.5
$ws.Range("M66").Value = -42418
$ws.Range("N66").Value = -14301.5
$ws.Range("H74").Value = 3998.125
$ws.Range("I74").Value = 4508.75
$ws.Range("J74").Value = 3487.5
$ws.Range("K74").Value = 4508.75
$ws.Range("L74").Value = 3487.5
$ws.Range("M74").Value = -3634.75
$ws.Range("N74").Value = -5235.5
$ws.Range("H77").Value = 3998.125
$ws.Range("I77").Value = 4508.75
$ws.Range("J77").Value = 3487.5
$ws.Range("K77").Value = 22543.75
$ws.Range("L77").Value = 17437.5
$ws.Range("M77").Value = -18175.75
$ws.Range("N77").Value = -26173.5
$ws.Range("H102").Value = 33335410
$ws.Range("I102").Value = 2345.625
$ws.Range("K102").Value = 2345.625
$ws.Range("M102").Value = -723.625
$ws.Range("H122").Value = 3285.724
$ws.Range("I122").Value = 2223.1177
$ws.Range("J122").Value = 4791.0835
$ws.Range("K122").Value = 6669.353099999999
$ws.Range("L122").Value = 14373.2505
$ws.Range("M122").Value = -4219.353099999999
$ws.Range("N122").Value = -19273.2505
$ws.Range("H132").Value = 2709.1292
$ws.Range("I132").Value = 1748.7368
$ws.Range("J132").Value = 4229.75
$ws.Range("K132").Value = 5246.2104
$ws.Range("L132").Value = 12689.25
$ws.Range("M132").Value = -2716.2104
$ws.Range("N132").Value = -17749.25
$ws.Range("H136").Value = 5867.3335
$ws.Range("I136").Value = 4344
$ws.Range("K136").Value = 13032
$ws.Range("M136").Value = -10482
$ws.Range("H139").Value = 80000.09
$ws.Range("J139").Value = 80000.10000000001
$ws.Range("L139").Value = 80000.10000000001
$ws.Range("N139").Value = -90280.10000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 131.875
$ws.Range("I7").Value = 108
$ws.Range("J7").Value = 203.5
$ws.Range("K7").Value = 108
$ws.Range("L7").Value = 203.5
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = -429.5
$ws.Range("H134").Value = 2581.9814
$ws.Range("I134").Value = 1662.4773
$ws.Range("K134").Value = 4987.4319
$ws.Range("M134").Value = -2452.4319

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2703.1667
$ws.Range("I31").Value = 1830.7878
$ws.Range("J31").Value = 4622.4
$ws.Range("K31").Value = 1830.7878
$ws.Range("L31").Value = 4622.4
$ws.Range("M31").Value = -1535.7878
$ws.Range("N31").Value = -5212.4
$ws.Range("H34").Value = 2703.1667
$ws.Range("I34").Value = 1830.7878
$ws.Range("J34").Value = 4622.4
$ws.Range("K34").Value = 1830.7878
$ws.Range("L34").Value = 4622.4
$ws.Range("M34").Value = -1628.7878
$ws.Range("N34").Value = -5026.4
$ws.Range("H58").Value = 6121.3335
$ws.Range("I58").Value = 4436.2
$ws.Range("K58").Value = 4436.2
$ws.Range("M58").Value = -4233.2
$ws.Range("H68").Value = 40000
$ws.Range("I68").Value = 40000
$ws.Range("K68").Value = 40000
$ws.Range("M68").Value = -39251
$ws.Range("H71").Value = 40000
$ws.Range("I71").Value = 40000
$ws.Range("K71").Value = 120000
$ws.Range("M71").Value = -116256
$ws.Range("H99").Value = 9762790
$ws.Range("I99").Value = 2444125
$ws.Range("K99").Value = 2444125
$ws.Range("M99").Value = -2442627
$ws.Range("H107").Value = 6298.222
$ws.Range("I107").Value = 790.9375
$ws.Range("K107").Value = 790.9375
$ws.Range("M107").Value = 1129.0625
$ws.Range("H126").Value = 9762790
$ws.Range("I126").Value = 2444125
$ws.Range("K126").Value = 7332375
$ws.Range("M126").Value = -7329905
$ws.Range("H132").Value = 2959.2222
$ws.Range("I132").Value = 1951.1333
$ws.Range("K132").Value = 5853.3999
$ws.Range("M132").Value = -3323.3999
$ws.Range("H136").Value = 6121.3335
$ws.Range("I136").Value = 4436.2
$ws.Range("K136").Value = 13308.6
$ws.Range("M136").Value = -10758.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1662.9
$ws.Range("J132").Value = 1644.4117
$ws.Range("L132").Value = 14799.7053
$ws.Range("N132").Value = -19859.7053

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 5848004.5
$ws.Range("J24").Value = 22006.625
$ws.Range("L24").Value = 22006.625
$ws.Range("N24").Value = -22352.625
$ws.Range("H80").Value = 20906584
$ws.Range("I80").Value = 189368.17
$ws.Range("J80").Value = 33336912
$ws.Range("K80").Value = 189368.17
$ws.Range("L80").Value = 33336912
$ws.Range("M80").Value = -188370.17
$ws.Range("N80").Value = -33338908
$ws.Range("H83").Value = 20906584
$ws.Range("I83").Value = 189368.17
$ws.Range("J83").Value = 33336912
$ws.Range("K83").Value = 946840.8500000001
$ws.Range("L83").Value = 166684560
$ws.Range("M83").Value = -941848.8500000001
$ws.Range("N83").Value = -166694544
$ws.Range("H126").Value = 5399.857
$ws.Range("J126").Value = 5966.5
$ws.Range("L126").Value = 17899.5
$ws.Range("N126").Value = -22839.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 220645.12
$ws.Range("I68").Value = 135913
$ws.Range("K68").Value = 135913
$ws.Range("M68").Value = -135164
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H71").Value = 220645.12
$ws.Range("I71").Value = 135913
$ws.Range("K71").Value = 679565
$ws.Range("M71").Value = -675821
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H98").Value = 79999
$ws.Range("J98").Value = 79999
$ws.Range("L98").Value = 79999
$ws.Range("N98").Value = -85989
$ws.Range("H122").Value = 5109.516
$ws.Range("I122").Value = 2666.7856
$ws.Range("K122").Value = 8000.3568
$ws.Range("M122").Value = -5550.3568
$ws.Range("H132").Value = 4430.86
$ws.Range("I132").Value = 3143.1538
$ws.Range("J132").Value = 8996.362999999999
$ws.Range("K132").Value = 9429.4614
$ws.Range("L132").Value = 26989.089
$ws.Range("M132").Value = -6899.4614
$ws.Range("N132").Value = -32049.089

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 316629.1
$ws.Range("I122").Value = 479466.75
$ws.Range("J122").Value = 5757.1816
$ws.Range("K122").Value = 1438400.25
$ws.Range("L122").Value = 17271.5448
$ws.Range("M122").Value = -1435950.25
$ws.Range("N122").Value = -22171.5448
